$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29; existing rows 29:95 shift down to 30:96
$ws.Rows("29:29").Insert()

# Populate the newly inserted row 29 with the new record
$ws.Cells.Item(29, 1).Value = 1
$ws.Cells.Item(29, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(29, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(29, 4).Value = 45203
$ws.Cells.Item(29, 5).Value = 15
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100103
$ws.Cells.Item(29, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(29, 9).Value = 100103004
$ws.Cells.Item(29, 10).Value = "Durazno"
$ws.Cells.Item(29, 11).Value = "Florida King"
$ws.Cells.Item(29, 12).Value = "Segunda"
$ws.Cells.Item(29, 13).Value = 250
$ws.Cells.Item(29, 14).Value = 28000
$ws.Cells.Item(29, 15).Value = 30000
$ws.Cells.Item(29, 16).Value = 28800
$ws.Cells.Item(29, 17).Value = "`$/bandeja 10 kilos granel"
$ws.Cells.Item(29, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(29, 19).Value = 2880
$ws.Cells.Item(29, 20).Value = 10
